# Update the dSF (column F) values for the listed rows as per the repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -2
    3  = 1
    4  = -4
    5  = -3
    6  = -2
    7  = -1
    10 = 6
    11 = -5
    12 = -1
    13 = -3
    14 = -2
    15 = -4
    16 = -1
    17 = -5
    18 = -4
    19 = 4
    20 = -1
    21 = 2
    22 = -3
    23 = 2
    24 = 3
    25 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
